# Update the "想去人数" (people interested) counts for two events.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows,
# so apply the same change to F5/F6 on each of them.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F5").Value = 7635
    $ws.Range("F6").Value = 5554
}
